$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 07:50"

# --- Update Estados Unidos row (row 4) ---
$ws.Range("B4").Value = 123774
$ws.Range("C4").Value = 196
$ws.Range("E4").Value = 118315
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 2228

# --- Update India row (row 44) ---
$ws.Range("D44").Value = 87
$ws.Range("E44").Value = 875
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 25

# --- Add Hungria as a new, updated entry right above Armenia (row 67) ---
# Insert a fresh row at 67 (old row 67 "Armenia" and everything below shifts down by one)
$ws.Rows("67").Insert()

$ws.Range("A67").Value = "Hungria"
$ws.Range("B67").Value = 408
$ws.Range("C67").Value = 65
$ws.Range("D67").Value = 34
$ws.Range("E67").Value = 361
$ws.Range("F67").Value = 6
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 13

# The old Hungria row (previously row 71, now shifted down to row 72 by the insert
# above) is now stale/duplicate data and must be removed; everything below shifts
# back up by one.
$ws.Rows("72").Delete()
